# Update confirmed/Indian case counts (and a few related columns) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 8    # Andhra Pradesh
$ws.Range("C5").Value = 3    # Bihar

$ws.Range("C8").Value = 32   # Gujarat
$ws.Range("D8").Value = 1

$ws.Range("C9").Value = 14   # Haryana

$ws.Range("C12").Value = 37  # Karnataka
$ws.Range("E12").Value = 3

$ws.Range("C13").Value = 87  # Kerala
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = 4

$ws.Range("C14").Value = 7   # Madhya Pradesh

$ws.Range("C15").Value = 86  # Maharashtra

$ws.Range("C16").Value = 1   # Manipur

$ws.Range("C21").Value = 29  # Punjab

$ws.Range("C22").Value = 30  # Rajasthan

$ws.Range("C24").Value = 13  # Tamil Nadu

$ws.Range("D27").Value = 1   # Uttarakhand

$ws.Range("C28").Value = 32  # Uttar Pradesh
$ws.Range("E28").Value = 11

$ws.Range("C29").Value = 9   # West Bengal

$ws.Range("C33").Value = 29  # Delhi
$ws.Range("E33").Value = 6
$ws.Range("F33").Value = 2
